$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Munka1")

# Row 9: distance-filter endpoint gets an extra :distance parameter
$ws.Range("A9").Value = "service/:lat/:lon/:distance"
$ws.Range("E9").Value = "szélesség, hosszúság, távolság"

# Row 14: keep same text (values unchanged, only internal shared-string
# reordering happened upstream), re-assert to be safe
$ws.Range("A14").Value = "service/my-services"
$ws.Range("C14").Value = "A felhasználó szolgáltatásainak lekérése"

# Update selection to match the saved view state
$ws.Range("E9").Select()
